$wb = $excel.ActiveWorkbook

# --- Sheet: Cardiac ---
$ws = $wb.Worksheets.Item("Cardiac")

# Row 2 and Row 3 swap their "info" text; the "Pain not worse with exertion" answer
# moves to row 3 and becomes TRUE, while "Do you have any PMHx?" moves to row 2 and
# is left unanswered (blank).
$ws.Range("B2").Value = "Do you have any PMHx? (counts as 2 independent minor features)"
$ws.Range("A2").Value = ""
$ws.Range("B3").Value = "Pain not worse with exertion (requires they clarify exercise 1hr after meal)"
$ws.Range("A3").Value = $true

# no associated shortness of breath -> TRUE
$ws.Range("A5").Value = $true

# no radiation to the neck, arm, or jaw? -> TRUE (was FALSE)
$ws.Range("A6").Value = $true

# Fix missing space typo in info text
$ws.Range("B9").Value = "Alternative cause of esoph dysphagia becomes obvious(food gets stuck or relieved by regurgitation of food)"

# no prior CAD -> TRUE
$ws.Range("A10").Value = $true

# no prior MI -> TRUE
$ws.Range("A13").Value = $true

# --- Sheet: GERD ---
$ws = $wb.Worksheets.Item("GERD")

# Alternative cause becomes obvious: esoph dysphagia -> blank (was TRUE)
$ws.Range("A6").Value = ""

# No dry cough -> blank (was TRUE)
$ws.Range("A10").Value = ""

# --- Sheet: Esophageal Dysphagia ---
$ws = $wb.Worksheets.Item("Esophageal Dysphagia")

# Food gets stuck -> TRUE (was blank)
$ws.Range("A2").Value = $true

# --- Sheet: CREST ---
$ws = $wb.Worksheets.Item("CREST")

# Difficulty swallowing liquids -> TRUE (was blank)
$ws.Range("A8").Value = $true

# no associated shortness of breath -> TRUE (was blank)
$ws.Range("A11").Value = $true
